# Automatic update of files.
# The source data rows for the "Artfynd" sheet got refreshed; the per-row
# observation data (id, sort order, red-list status, taxon id, names,
# author, locality, and coordinates) moved between rows 4, 5, 6, 8 and 9:
#   new row 4 <- old row 9
#   new row 5 <- old row 4
#   new row 6 <- old row 8
#   new row 8 <- old row 6
#   new row 9 <- old row 5
# Rows 3 and 7 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")
$rows = @(4, 5, 6, 8, 9)

# Snapshot current values for every relevant cell before we overwrite anything.
# NOTE: read via .Value2 (not the plain .Value indexed property, which this
# COM shim cannot resolve without an explicit index and returns a useless
# reflection description for instead).
$snapshot = @{}
foreach ($r in $rows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# new row -> source (old) row
$mapping = @{ 4 = 9; 5 = 4; 6 = 8; 8 = 6; 9 = 5 }

foreach ($dstRow in $rows) {
    $srcRow = $mapping[$dstRow]
    foreach ($c in $cols) {
        $dstAddr = "$c$dstRow"
        $srcAddr = "$c$srcRow"
        $ws.Range($dstAddr).Value = $snapshot[$srcAddr]
    }
}
